# Scheduled runner update: refresh market-price / profit figures on the
# Anima_Profits sheets (one tab per crafting job) with the latest values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6067.143
$ws.Range("I28").Value = 216.26666
$ws.Range("J28").Value = 20694.334
$ws.Range("K28").Value = 216.26666
$ws.Range("L28").Value = 20694.334
$ws.Range("M28").Value = 268.73334
$ws.Range("N28").Value = -21664.334

$ws.Range("H41").Value = 229.6
$ws.Range("I41").Value = 170.85715
$ws.Range("J41").Value = 366.66666
$ws.Range("K41").Value = 170.85715
$ws.Range("L41").Value = 366.66666
$ws.Range("M41").Value = 269.14285
$ws.Range("N41").Value = -1246.66666

$ws.Range("H76").Value = 5627.273
$ws.Range("I76").Value = 7000
$ws.Range("J76").Value = 3225
$ws.Range("K76").Value = 7000
$ws.Range("L76").Value = 3225
$ws.Range("M76").Value = -6685
$ws.Range("N76").Value = -3855

$ws.Range("H79").Value = 5627.273
$ws.Range("I79").Value = 7000
$ws.Range("J79").Value = 3225
$ws.Range("K79").Value = 7000
$ws.Range("L79").Value = 3225
$ws.Range("M79").Value = -5908
$ws.Range("N79").Value = -5409

$ws.Range("H137").Value = 1263.12
$ws.Range("I137").Value = 1169.1177
$ws.Range("J137").Value = 1462.875
$ws.Range("K137").Value = 3507.3531
$ws.Range("L137").Value = 4388.625
$ws.Range("M137").Value = -957.3531000000003
$ws.Range("N137").Value = -9488.625

$ws.Range("H138").Value = 2667.6924
$ws.Range("I138").Value = 3269.75
$ws.Range("J138").Value = 2487.075
$ws.Range("K138").Value = 9809.25
$ws.Range("L138").Value = 7461.224999999999
$ws.Range("M138").Value = -4669.25
$ws.Range("N138").Value = -17741.225

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 10878.692
$ws.Range("I37").Value = 2001
$ws.Range("J37").Value = 11618.5
$ws.Range("K37").Value = 2001
$ws.Range("L37").Value = 11618.5
$ws.Range("M37").Value = -1728
$ws.Range("N37").Value = -12164.5

$ws.Range("H88").Value = 2834
$ws.Range("I88").Value = 2778.6667
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 2778.6667
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -2372.6667
$ws.Range("N88").Value = -3812

$ws.Range("H91").Value = 2834
$ws.Range("I91").Value = 2778.6667
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 2778.6667
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -1374.6667
$ws.Range("N91").Value = -5808

$ws.Range("H132").Value = 4182.838
$ws.Range("I132").Value = 4441.643
$ws.Range("J132").Value = 3377.6667
$ws.Range("K132").Value = 13324.929
$ws.Range("L132").Value = 10133.0001
$ws.Range("M132").Value = -10794.929
$ws.Range("N132").Value = -15193.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 11903.5
$ws.Range("I82").Value = 2327.1
$ws.Range("J82").Value = 23874
$ws.Range("K82").Value = 2327.1
$ws.Range("L82").Value = 23874
$ws.Range("M82").Value = -1944.1
$ws.Range("N82").Value = -24640

$ws.Range("H85").Value = 11903.5
$ws.Range("I85").Value = 2327.1
$ws.Range("J85").Value = 23874
$ws.Range("K85").Value = 2327.1
$ws.Range("L85").Value = 23874
$ws.Range("M85").Value = -1001.1
$ws.Range("N85").Value = -26526

$ws.Range("H86").Value = 4600
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -7246

$ws.Range("H89").Value = 4600
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -36232

$ws.Range("H134").Value = 3066.0967
$ws.Range("I134").Value = 3056.1667
$ws.Range("J134").Value = 3079.8462
$ws.Range("K134").Value = 9168.500100000001
$ws.Range("L134").Value = 9239.5386
$ws.Range("M134").Value = -6633.500100000001
$ws.Range("N134").Value = -14309.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 18365.834
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 18365.834
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 18365.834
$ws.Range("N51").Value = -19837.834

$ws.Range("H59").Value = 26410.572
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 26410.572
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 26410.572
$ws.Range("N59").Value = -28700.572

$ws.Range("H60").Value = 10667.167
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 10667.167
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 10667.167
$ws.Range("N60").Value = -11689.167

$ws.Range("H61").Value = 18365.834
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 18365.834
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 18365.834
$ws.Range("N61").Value = -19061.834

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0

$ws.Range("H68").Value = 1354.4
$ws.Range("I68").Value = 750.3929000000001
$ws.Range("J68").Value = 1589.2916
$ws.Range("K68").Value = 2251.1787
$ws.Range("L68").Value = 4767.8748
$ws.Range("M68").Value = -1440.1787
$ws.Range("N68").Value = -6389.8748

$ws.Range("H71").Value = 1354.4
$ws.Range("I71").Value = 750.3929000000001
$ws.Range("J71").Value = 1589.2916
$ws.Range("K71").Value = 6753.5361
$ws.Range("L71").Value = 14303.6244
$ws.Range("M71").Value = -2697.5361
$ws.Range("N71").Value = -22415.6244

$ws.Range("H98").Value = 333596.66
$ws.Range("I98").Value = 290
$ws.Range("J98").Value = 500250
$ws.Range("K98").Value = 870
$ws.Range("L98").Value = 1500750
$ws.Range("M98").Value = 628
$ws.Range("N98").Value = -1503746

$ws.Range("H113").Value = 1037.5312
$ws.Range("I113").Value = 577
$ws.Range("J113").Value = 1498.0625
$ws.Range("K113").Value = 1731
$ws.Range("L113").Value = 4494.1875
$ws.Range("M113").Value = 439
$ws.Range("N113").Value = -8834.1875

$ws.Range("H131").Value = 1035.2632
$ws.Range("I131").Value = 715
$ws.Range("J131").Value = 1183.0769
$ws.Range("K131").Value = 2145
$ws.Range("L131").Value = 3549.2307
$ws.Range("M131").Value = 2895
$ws.Range("N131").Value = -13629.2307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 2290
$ws.Range("I55").Value = 985
$ws.Range("J55").Value = 4900
$ws.Range("K55").Value = 985
$ws.Range("L55").Value = 4900
$ws.Range("M55").Value = -658
$ws.Range("N55").Value = -5554

$ws.Range("H70").Value = 5888.8213
$ws.Range("I70").Value = 5743.9
$ws.Range("J70").Value = 6251.125
$ws.Range("K70").Value = 5743.9
$ws.Range("L70").Value = 6251.125
$ws.Range("M70").Value = -5473.9
$ws.Range("N70").Value = -6791.125

$ws.Range("H73").Value = 5888.8213
$ws.Range("I73").Value = 5743.9
$ws.Range("J73").Value = 6251.125
$ws.Range("K73").Value = 5743.9
$ws.Range("L73").Value = 6251.125
$ws.Range("M73").Value = -4807.9
$ws.Range("N73").Value = -8123.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2496.608
$ws.Range("I132").Value = 2042.9722
$ws.Range("J132").Value = 3585.3333
$ws.Range("K132").Value = 6128.9166
$ws.Range("L132").Value = 10755.9999
$ws.Range("M132").Value = -3598.9166
$ws.Range("N132").Value = -15815.9999
